$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before A; this shifts the existing columns (and their
# formatting / column widths) one place to the right, turning the old
# A..I (last_name..email) into B..J.
$ws.Columns("A").Insert()

# Populate the new "first_name" column.
$ws.Range("A1").Value = "first_name"
$ws.Range("A2").Value = "John"

# The old "last_name" column (now column B) keeps "last_name" as its header
# but the sample value becomes "Doe" (since "John" moved to the new column).
$ws.Range("B2").Value = "Doe"

# Match the formatting of the new column to the "target_degree" column
# (general alignment, no border/font override) like the rest of the sheet.
$ws.Range("A1:A2").HorizontalAlignment = 1
